# Material.xlsx update: add "Multimedia Folder" column (I) used by the
# nodegoat scripts, now renamed to process_data, and mark each row's
# processing stage ("output" for the raw data-import row, "processed"
# for every csv-export row). Also narrow column F now that the old long
# bestFit text no longer needs the extra width.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column I: header + values ------------------------------------
$ws.Range("I1").Value = "Multimedia Folder"
$ws.Range("I2").Value = "output"
$ws.Range("I3:I20").Value = "processed"

# Header I1 should carry the same bold/bordered header style as the rest
# of row 1 (copy format from H1, which already has it).
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Column widths ------------------------------------------------------
# Column F no longer needs to auto-fit its old (much longer) text, so it
# gets a fixed, narrower width; new column I gets a fixed width too.
$ws.Columns.Item(6).ColumnWidth = 17.5
$ws.Columns.Item(9).ColumnWidth = 16.666666666666668

# --- Selection ------------------------------------------------------
$ws.Range("I3:I20").Select() | Out-Null
